# Third commit including jenkins
#
# The LoginData sheet gains a new C2 cell recording the most recent
# "Logged in as ..." confirmation message captured by the Jenkins test run.
# (LoginData is not the active tab - CheckOut is - so address it by name.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

$ws.Range("C2").Value = "Logged in as standard_user_20241227_010034"
